$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values updated for rows 2-5
$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 3
